# Append seven new English/Chinese vocabulary pairs to the word list on
# Sheet1 (rows 175-181, columns A/B) and move the active selection to
# D177, matching the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @("replace",      "代替"),
    @("background",   "背景"),
    @("feature",      "特征"),
    @("enhancements", "增强功能"),
    @("detail",       "细节"),
    @("condense",     "压缩"),
    @("majority",     "大多数")
)

$row = 175
foreach ($pair in $pairs) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Range("D177").Select()
